$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-23 Thursday" "2025-10-24 Friday"

Replace-Text "703÷5=140, 3" "672÷8=84, 0"
Replace-Text "219÷3=73, 0" "885÷3=295, 0"
Replace-Text "885÷5=177, 0" "759÷6=126, 3"
Replace-Text "109÷8=13, 5" "952÷7=136, 0"
Replace-Text "686÷3=228, 2" "961÷2=480, 1"

Replace-Text "941÷2=470, 1" "449÷4=112, 1"
Replace-Text "637÷7=91, 0" "561÷5=112, 1"
Replace-Text "122÷6=20, 2" "893÷9=99, 2"
Replace-Text "808÷7=115, 3" "937÷5=187, 2"
Replace-Text "900÷6=150, 0" "661÷3=220, 1"

Replace-Text "762÷9=84, 6" "676÷3=225, 1"
Replace-Text "499÷4=124, 3" "660÷5=132, 0"
Replace-Text "501÷7=71, 4" "436÷8=54, 4"
Replace-Text "309÷4=77, 1" "931÷4=232, 3"
Replace-Text "361÷8=45, 1" "984÷7=140, 4"

Replace-Text "847÷3=282, 1" "288÷6=48, 0"
Replace-Text "860÷9=95, 5" "873÷8=109, 1"
Replace-Text "398÷7=56, 6" "965÷2=482, 1"
Replace-Text "742÷4=185, 2" "101÷6=16, 5"
Replace-Text "869÷2=434, 1" "334÷6=55, 4"

Replace-Text "190÷6=31, 4" "457÷3=152, 1"
Replace-Text "228÷8=28, 4" "245÷5=49, 0"
Replace-Text "734÷2=367, 0" "353÷4=88, 1"
Replace-Text "591÷9=65, 6" "793÷3=264, 1"
Replace-Text "239÷8=29, 7" "452÷7=64, 4"

Write-Host "Done."
